$wb = $excel.ActiveWorkbook

$oldGuid = "4db34cc0-5ced-499a-accb-c6f86f197e33"
$newGuid = "1a99bda7-55d2-4d90-9662-5ba4908eab92"

$oldHash = "2ed6759a3c6885a1f2959f87875afc7833c0a8ae"
$newHash = "66861ffe82cf035077894158410bd0cbf7567fa9"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# Overview sheet: File Name / Path And Name / Latest HO Xliff Generate Date
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-31 11:03:50"

# zh-cn sheet: Source File Name / Latest Handoff File / Latest Handoff Datetime
$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-31 11:03:46"

# de-de sheet: Source File Name / Latest Handoff File
$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
# de-de!H2 ("Latest Handoff Datetime") shares the same source string as
# Overview!G2 ("2016-08-31 11:03:32" -> "...:50"), so it moves too.
$wsDe.Range("H2").Value = "2016-08-31 11:03:50"

# Update hyperlink display text to match the new file names
# (iterating the worksheet's Hyperlinks collection is required for the
#  property write to land on the existing hyperlink instead of minting a
#  stray duplicate)
foreach ($h in $wsOverview.Hyperlinks) {
    if ($h.TextToDisplay -eq "e2e\$oldGuid.md") {
        $h.TextToDisplay = "e2e\$newGuid.md"
    }
}

foreach ($h in $wsZh.Hyperlinks) {
    if ($h.TextToDisplay -eq "$oldGuid.md") {
        $h.TextToDisplay = "$newGuid.md"
    }
}

foreach ($h in $wsDe.Hyperlinks) {
    if ($h.TextToDisplay -eq "$oldGuid.md") {
        $h.TextToDisplay = "$newGuid.md"
    }
}
